# feat: hide name plate option
#
# Inserts a new row for the "hideNamePlateButton" setting into the "Text"
# sheet (row 36 — just after the "Hide Player Names" block), shifting every
# row below it down by one, and makes the "Text" sheet the active tab with
# the selection parked on the last data row (mirrors where Excel would leave
# the cursor after typing the new row and scrolling to the bottom).

$wb = $excel.ActiveWorkbook
$wsText = $wb.Worksheets.Item("Text")

# Insert a new blank row at 36, pushing rows 37.. down to 38..
$wsText.Rows.Item(36).Insert()

# English key / English display text / Japanese display text for the new
# "Hide Nameplate" option.
$wsText.Range("A36").Value = "hideNamePlateButton"
$wsText.Range("B36").Value = "Hide Nameplate"
$wsText.Range("M36").Value = "ネームプレートを非表示にする"

# "Text" becomes the selected/active sheet (was "Neutral").
$wsText.Activate()
$wsText.Range("M127").Select() | Out-Null
